$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stickers")

# Add new row 4 first (reuses the shared-string slots freed up once B2/B3 change below)
$ws.Range("A4").Value = "как дела?"
$ws.Range("B4").Value = "CAACAgIAAxkBAANZYhKp4H4SJJwrTOqu6UEiq9EtqZ0AAhMBAAJSiZEjgGq_p-zb_zwjBA"

# Update existing rows 2 and 3 (B column) with the new sticker ids
$ws.Range("B2").Value = "CAACAgIAAxkBAANhYhKrDYXggOZJnoRRnorQLIi9TG4AAh4JAAIYQu4I-VjZ7h0hnCEjBA"
$ws.Range("B3").Value = "CAACAgIAAxkBAANeYhKq-3BVuhgAAXmf_WK95nib4jAFAAKOAAMWQmsKvqSGfW1-LVwjBA"

$ws.Range("A5").Value = "удачи"
$ws.Range("B5").Value = "CAACAgIAAxkBAANkYhKrPLhmoae0e9K8m5jF2H45Wn0AAs8AA_cCyA-5-Dj7pxiu_SME"

$ws.Range("A6").Value = "знаешь?"
$ws.Range("B6").Value = "CAACAgIAAxkBAANnYhKrQE9zIKlx04LT4nLcNTCrMIYAAs4AA_cCyA9cmEfgzKtuiiME"

# Copy the style used on B2/B3 down to the new B4:B6 cells
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B4:B6").PasteSpecial(-4122) | Out-Null

# Update the active selection to match the saved view state
$ws.Range("L8").Select() | Out-Null
